$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.453.61'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '1.795.42'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.99'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.555'
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.47'
$ws.Range('E8').Value = '  +1.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.296'
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0693'
$ws.Range('E10').Value = '  +0.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '2.053.31'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.10'
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.796.57'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.633'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '34.398.26'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.24'
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.35'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '0.0₃0801'
$ws.Range('E19').Value = '  +2.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '246.50'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.09'
$ws.Range('E21').Value = '  +1.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.16'
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.71'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.25'
$ws.Range('E26').Value = '  +1.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.42'
$ws.Range('E27').Value = '  +0.54%  '
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +0.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0523'
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.90'
$ws.Range('E32').Value = '  +8.25%  '
$ws.Range('E33').Value = '  +3.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.83'
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('D35').Value = '1.443.61'
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.61'
$ws.Range('E36').Value = '  +6.74%  '
$ws.Range('E37').Value = '  +4.00%  '
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('E39').Value = '  +1.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '84.08'
$ws.Range('E40').Value = '  +4.63%  '
$ws.Range('E41').Value = '  +1.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.933'
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.75'
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.81'
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('E45').Value = '  +3.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.09'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').Value = '1.949.60'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.69'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = '0.0₆0129'
